$d = $word.ActiveDocument

# Locate the "Testing" top-level bullet paragraph (it sits right before
# "Results and conclusion") so the three new sub-bullets can be inserted
# directly after it.
$testingIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd("`r`a") -eq "Testing") {
        $testingIndex = $i
    }
}

$newItems = @("Testing group", "Questionnaire", "results")

$insertIndex = $testingIndex
foreach ($text in $newItems) {
    $insertIndex = $insertIndex + 1

    # Create a new paragraph right after the previous one.
    $prevPara = $d.Paragraphs.Item($insertIndex - 1)
    $prevPara.Range.InsertParagraphAfter()

    # The freshly created paragraph inherits the "ListParagraph" style and
    # numId from "Testing", but at the same outline level (ilvl 0). Bump it
    # down to ilvl 1 (ListLevelNumber is 1-based, so level 1 == value 2),
    # then fill in its text.
    $newPara = $d.Paragraphs.Item($insertIndex)
    $newPara.Range.ListFormat.ListLevelNumber = 2
    $newPara.Range.InsertAfter($text)
}
